$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '37.021.50'
Set-TextValue 'E2' '  -1.19%  '

Set-TextValue 'D3' '2.011.90'
Set-TextValue 'E3' '  -1.89%  '

Set-TextValue 'E4' '  +0.03%  '

Set-TextValue 'D5' '226.35'
Set-TextValue 'E5' '  -1.00%  '

Set-TextValue 'D6' '0.605'
Set-TextValue 'E6' '  -1.06%  '

Set-TextValue 'E7' '  -0.01%  '

Set-TextValue 'D8' '55.20'
Set-TextValue 'E8' '  -2.02%  '

Set-TextValue 'D9' '0.373'
Set-TextValue 'E9' '  -3.43%  '

Set-TextValue 'D10' '0.0776'
Set-TextValue 'E10' '  -3.75%  '

Set-TextValue 'E11' '  -4.15%  '

Set-TextValue 'D12' '2.309.90'
Set-TextValue 'E12' '  -1.89%  '

Set-TextValue 'D13' '14.04'
Set-TextValue 'E13' '  -3.51%  '

Set-TextValue 'D14' '19.77'
Set-TextValue 'E14' '  -4.17%  '

Set-TextValue 'B15' 'Polygon'
Set-TextValue 'C15' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D15' '0.735'
Set-TextValue 'E15' '  -2.53%  '

Set-TextValue 'B16' 'Polkadot'
Set-TextValue 'C16' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D16' '5.16'
Set-TextValue 'E16' '  -1.91%  '

Set-TextValue 'D17' '2.013.25'
Set-TextValue 'E17' '  -2.88%  '

Set-TextValue 'D18' '36.971.45'
Set-TextValue 'E18' '  -0.97%  '

Set-TextValue 'D19' '6.21'
Set-TextValue 'E19' '  +1.94%  '

Set-TextValue 'D20' '68.78'
Set-TextValue 'E20' '  -1.52%  '

Set-TextValue 'D21' '0.0₃0810'
Set-TextValue 'E21' '  -3.98%  '

Set-TextValue 'D22' '222.73'
Set-TextValue 'E22' '  -1.43%  '

Set-TextValue 'E23' '  -0.01%  '

Set-TextValue 'E24' '  +2.15%  '

Set-TextValue 'D25' '2.18'
Set-TextValue 'E25' '  -4.22%  '

Set-TextValue 'D26' '164.48'
Set-TextValue 'E26' '  -2.35%  '

Set-TextValue 'D27' '8.95'
Set-TextValue 'E27' '  -5.99%  '

Set-TextValue 'D28' '0.126'
Set-TextValue 'E28' '  -2.43%  '

Set-TextValue 'D29' '18.57'
Set-TextValue 'E29' '  -1.88%  '

Set-TextValue 'D30' '1.31'
Set-TextValue 'E30' '  -4.93%  '

Set-TextValue 'D31' '0.117'
Set-TextValue 'E31' '  -1.18%  '

Set-TextValue 'D32' '4.40'
Set-TextValue 'E32' '  -3.29%  '

Set-TextValue 'D33' '0.0599'
Set-TextValue 'E33' '  -2.21%  '

Set-TextValue 'D34' '4.44'
Set-TextValue 'E34' '  -2.47%  '

Set-TextValue 'D35' '2.33'
Set-TextValue 'E35' '  -2.81%  '

Set-TextValue 'D36' '1.86'
Set-TextValue 'E36' '  +2.21%  '

Set-TextValue 'E37' '  +0.18%  '

Set-TextValue 'D38' '3.14'
Set-TextValue 'E38' '  -1.78%  '

Set-TextValue 'D39' '5.38'
Set-TextValue 'E39' '  -0.89%  '

Set-TextValue 'D40' '1.464.11'
Set-TextValue 'E40' '  -2.60%  '

Set-TextValue 'D41' '0.0211'
Set-TextValue 'E41' '  -4.41%  '

Set-TextValue 'D42' '94.26'
Set-TextValue 'E42' '  -2.01%  '

Set-TextValue 'B43' 'Cronos'
Set-TextValue 'C43' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D43' '0.0912'
Set-TextValue 'E43' '  -2.31%  '

Set-TextValue 'B44' 'HuobiToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D44' '2.76'
Set-TextValue 'E44' '  -4.01%  '

Set-TextValue 'B45' 'FTXToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 'D45' '4.17'
Set-TextValue 'E45' '  +9.86%  '

Set-TextValue 'D46' '1.12'
Set-TextValue 'E46' '  -2.47%  '

Set-TextValue 'D47' '15.87'
Set-TextValue 'E47' '  -6.13%  '

Set-TextValue 'D48' '0.998'
Set-TextValue 'E48' '  -1.91%  '

Set-TextValue 'D49' '7.05'
Set-TextValue 'E49' '  -2.30%  '

Set-TextValue 'E50' '  -1.08%  '

Set-TextValue 'D51' '2.196.89'
Set-TextValue 'E51' '  -1.86%  '
